# Insert a new "latest price snapshot" column just before the "nom" column.
# Before: ... BG=<last price>, BH=nom, BI=url_produit
# After:  ... BG=<last price>, BH=<new price snapshot>, BI=nom, BJ=url_produit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BH ("nom") is column 60. Inserting here shifts BH->BI and BI->BJ,
# carrying over formatting (so the bold/bordered header style on row 1
# moves along with the old BH/BI cells) exactly like Excel's own
# "Insert Column" command.
$ws.Columns("BH").Insert()

# New header cell for the freshly inserted price-history column.
$ws.Range("BH1").Value = "2026-01-30 10:22:49"

# Figure out the last used row from the sheet dimensions (206 data rows
# plus the header row in this workbook, but compute it dynamically so the
# script keeps working if the sheet grows).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# For every product row, the new column simply repeats the most recent
# known price (the value that was already sitting in column BG, the last
# timestamped price column) - matching rows that have no recorded price
# simply stay blank.
for ($r = 2; $r -le $lastRow; $r++) {
    $v = $ws.Cells.Item($r, 59).Value2
    if ($v -ne $null -and $v -ne "") {
        $ws.Cells.Item($r, 60).Value2 = $v
    }
}
